$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B4 text: remove trailing space from "Separator "
$ws.Range("B4").Value = "Separator"

# Apply bold Aptos 9.8pt font with vertical-centered alignment to B4
$ws.Range("B4").Font.Bold = $true
$ws.Range("B4").Font.Name = "Aptos"
$ws.Range("B4").Font.Size = 9.8
$ws.Range("B4").VerticalAlignment = -4108

# Update selection to E6
$ws.Range("E6").Select()
